$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The title paragraph ("On Pilgrimage - December 1978", styled Heading1) is
# wrapped in a bookmark whose name ("on-pilgrimage---december-1978") is not a
# legal Word bookmark identifier (it contains hyphens), so it is invisible to
# the Bookmarks collection. We get rid of it by deleting the whole first
# paragraph outright (which collapses both the bookmarkStart/bookmarkEnd
# markers down to the start of the document) and then stripping those two
# now-adjacent, zero-length markers explicitly before rebuilding the
# paragraph with the new pandoc-style title-block structure.
# ---------------------------------------------------------------------------

$titleParaEnd = $d.Paragraphs.Item(1).Range.End
$d.Range(0, $titleParaEnd).Delete()

# Remove the two bookmark markers that collapsed onto position 0.
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# ---------------------------------------------------------------------------
# Paragraph 1 is now "By Dorothy Day" (bold run). Turn it into an "Authors"
# styled paragraph containing "Dorothy Day", split word-by-word/space into
# separate runs (mirrors how pandoc emits its title-block author line).
# ---------------------------------------------------------------------------

$authorsXml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>'

$d.Paragraphs.Item(1).Range.InsertXML($authorsXml)

# ---------------------------------------------------------------------------
# Insert a new "Title" styled paragraph in front of it holding
# "On Pilgrimage - December 1978", again split word-by-word/space into
# separate runs.
# ---------------------------------------------------------------------------

$titleXml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">On</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">-</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">December</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">1978</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>'

$d.Range(0, 0).InsertXML($titleXml)

Write-Output ("p1 [" + $d.Paragraphs.Item(1).Range.Text + "] style=" + $d.Paragraphs.Item(1).Style.NameLocal)
Write-Output ("p2 [" + $d.Paragraphs.Item(2).Range.Text + "] style=" + $d.Paragraphs.Item(2).Style.NameLocal)
Write-Output ("total paragraphs: " + $d.Paragraphs.Count)
